# Update attendance/figure counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 194
$wsExhibit.Range("F4").Value = 2260
$wsExhibit.Range("F5").Value = 1716
$wsExhibit.Range("F6").Value = 327
$wsExhibit.Range("F8").Value = 784

# --- Sheet "全部类型" (all types, aggregated) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 194
$wsAll.Range("F4").Value = 2260
$wsAll.Range("F5").Value = 1716
$wsAll.Range("F6").Value = 327
$wsAll.Range("F9").Value = 783
